# Header row (row 1): A=Date, B=Roll, C=Name, D=Total Attendance Count,
#                      E=Real, F=Duplicate, G=Invalid, H=Absent
# Update attendance values for each dated row (rows 3-18) per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

$ws.Range("H5").Value = 1

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

$ws.Range("H7").Value = 1

$ws.Range("H8").Value = 1

$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

$ws.Range("H11").Value = 1

$ws.Range("H12").Value = 1

$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1

$ws.Range("H14").Value = 1

$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1

$ws.Range("H16").Value = 1

$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1

$ws.Range("H18").Value = 1
